$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 label becomes "List exp.: NR/GCS/C/O", values restored from former row 5)
$ws.Range("A2").Value = "List exp.: NR/GCS/C/O"
$ws.Range("B2").Value = 1.95572072394667
$ws.Range("C2").Value = 2.52780478947738
$ws.Range("D2").Value = 2.78980832908363
$ws.Range("E2").Value = 2.17878283022247
$ws.Range("F2").Value = 2.66429365853611
$ws.Range("G2").Value = 2.58649963757343

# Row 3 (A3 label becomes "List exp.: NR/C/O", values restored from former row 4)
$ws.Range("A3").Value = "List exp.: NR/C/O"
$ws.Range("B3").Value = 1.40117746977089
$ws.Range("C3").Value = 1.85219004920533
$ws.Range("D3").Value = 2.07905099352898
$ws.Range("E3").Value = 1.63238774908649
$ws.Range("F3").Value = 1.99676754389218
$ws.Range("G3").Value = 1.92594993505939

# Row 4 (A4 label becomes "List exp.: GCS/C/O", values restored from former row 3)
$ws.Range("A4").Value = "List exp.: GCS/C/O"
$ws.Range("B4").Value = 1.40098726783718
$ws.Range("C4").Value = 1.87398941103836
$ws.Range("D4").Value = 2.0708280327613
$ws.Range("E4").Value = 1.70749447477815
$ws.Range("F4").Value = 1.92760940318959
$ws.Range("G4").Value = 1.83955709804286

# Row 5 (A5 label becomes "List exp.: C/O", values restored from former row 2)
$ws.Range("A5").Value = "List exp.: C/O"
$ws.Range("B5").Value = 0.910926213838564
$ws.Range("C5").Value = 1.11363712713086
$ws.Range("D5").Value = 1.29297415821355
$ws.Range("E5").Value = 0.931581478142742
$ws.Range("F5").Value = 1.08083416372115
$ws.Range("G5").Value = 1.25453452227073
